$wb = $excel.ActiveWorkbook

# Sheets: About, SoCDTtiNTY-psgr, SoCDTtiNTY-frgt
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# Update LDVs row (row 2) on psgr sheet: 7.6% -> 7.55%
$wsPsgr.Range("B2:H2").Value = 0.0755

# Update HDVs row (row 3) on frgt sheet with new values
$wsFrgt.Range("B3").Value = 0.0219
$wsFrgt.Range("C3").Value = 0.0219
$wsFrgt.Range("D3").Value = 0.081
$wsFrgt.Range("E3").Value = 0.0245
$wsFrgt.Range("F3").Value = 0.0219
$wsFrgt.Range("G3").Value = 0.0219
$wsFrgt.Range("H3").Value = 0.0219

# Update selections
$wsPsgr.Range("E14").Select()
$wsFrgt.Range("B2:H7").Select()

# Make frgt sheet active (tab selected)
$wsFrgt.Activate()
